$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.492.02'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.914.01'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4833'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2893'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06712'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '110.09'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.86'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.13%  '
$ws.Range('D12').Value = '1.914.50'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07555'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.275'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6690'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '276.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '30.501.84'
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.0000'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').Value = '2.163.48'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.482'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.436'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.439'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.120'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1054'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.151'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.053'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04992'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7293'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.131'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9991'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.732'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02031'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.668'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '110.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.019'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4429'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8649'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('E44').Value = '  +0.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.363'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.265'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.06%  '

$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '47.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.42%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1237'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.85%  '
